$d = $word.ActiveDocument

$pairs = @(
    @("889÷5=", "161÷6="),
    @("162÷9=", "755÷4="),
    @("832÷8=", "758÷8="),
    @("429÷7=", "899÷4="),
    @("528÷3=", "706÷9="),
    @("190÷4=", "238÷7="),
    @("735÷5=", "484÷5="),
    @("841÷3=", "938÷4="),
    @("112÷9=", "208÷9="),
    @("815÷4=", "241÷6="),
    @("328÷9=", "256÷6="),
    @("810÷6=", "219÷8="),
    @("971÷3=", "573÷5="),
    @("176÷7=", "549÷9="),
    @("136÷5=", "367÷9="),
    @("670÷2=", "474÷6="),
    @("463÷2=", "141÷2="),
    @("143÷2=", "624÷7="),
    @("996÷5=", "965÷5="),
    @("861÷9=", "899÷7="),
    @("450÷6=", "427÷3="),
    @("829÷9=", "404÷5="),
    @("222÷3=", "582÷4="),
    @("968÷7=", "525÷3="),
    @("106÷7=", "166÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
